# The commit swaps the two theme parts of the deck: the slide-master
# theme (ppt/theme/theme1.xml, originally the "Integral" / "Red Violet"
# palette) becomes the stock "Office Theme" palette, while the notes
# master's theme (ppt/theme/theme2.xml) becomes the old "Integral" /
# "Red Violet" palette.
#
# The PowerPoint object model only exposes per-slot colour editing of
# the *active* (slide-master) theme via ThemeColorScheme.Colors(i).RGB
# (indices 1-12 => dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) -
# there is no Master/NotesMaster.ThemeColorScheme in the object model,
# so we push the "Office Theme" palette onto theme1.xml that way.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

# Office Theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# expressed as COM BGR-packed RGB integers.
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $cs.Item($i).RGB = $officeThemeColors[$i - 1]
}
